# Applies the "Refined metadata to be additional tab" change:
#  1. Updates the panel_query_time / time_taken timestamps on the "data" sheet.
#  2. Adds a new "metadata" worksheet (after "data") describing the panel query.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- 1. Refresh the time_taken timestamps on the data sheet (column F, rows 2-16) ---
$data.Range("F2").Value  = "2021-10-05 14:20:43.681253"
$data.Range("F3").Value  = "2021-10-05 14:20:43.681261"
$data.Range("F4").Value  = "2021-10-05 14:20:43.681264"
$data.Range("F5").Value  = "2021-10-05 14:20:43.681267"
$data.Range("F6").Value  = "2021-10-05 14:20:43.681270"
$data.Range("F7").Value  = "2021-10-05 14:20:43.681273"
$data.Range("F8").Value  = "2021-10-05 14:20:43.681276"
$data.Range("F9").Value  = "2021-10-05 14:20:43.681279"
$data.Range("F10").Value = "2021-10-05 14:20:43.681282"
$data.Range("F11").Value = "2021-10-05 14:20:43.681284"
$data.Range("F12").Value = "2021-10-05 14:20:43.681287"
$data.Range("F13").Value = "2021-10-05 14:20:43.681289"
$data.Range("F14").Value = "2021-10-05 14:20:43.681292"
$data.Range("F15").Value = "2021-10-05 14:20:43.681295"
$data.Range("F16").Value = "2021-10-05 14:20:43.681297"

# --- 2. Add the new "metadata" worksheet right after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Reuse the bold/bordered header style already used on the "data" sheet (B1:F1)
# for the new header row, and the header-row style used on A2 for the numeric
# index cell, by copying formats across instead of re-creating new styles.
$data.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)

$data.Range("B1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Hereditary Erythrocytosis"
$meta.Range("C2").Value = 157

# "1.35" must stay a text value (not get coerced to a number) but without
# picking up a quote-prefix / number-format style, so format as text, set the
# value, then clear the temporary text format back to the plain (unstyled) look.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.35"
$meta.Range("D2").ClearFormats()

$meta.Range("E2").Value = "2021-03-23T10:40:37.534764Z"
$meta.Range("F2").Value = "2021-10-05 14:20:43.677720"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/157/?format=json"

Write-Output "metadata sheet added; timestamps refreshed"
